$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tournament results entered below the existing data (rows 160-175).
# Columns: A=Player1, B=Player2, C=Player3, D=Teamname, E=Datum, F=Platzierung,
#          G=Anzahl Teams, H=Draft
$rows = @(
    @{A='Chris';    B='Valdemar'; C='';      D='Pferdewelle Stampf Stampf';   E=45479; F=3; G=4; H=0}
    @{A='Philipp';  B='André';    C='';      D='Garruk Ultras';               E=45479; F=1; G=4; H=0}
    @{A='Marie';    B='Emilio';   C='';      D='Die Prenzlauer Crew';         E=45479; F=3; G=4; H=0}
    @{A='Phia';     B='Daisy';    C='';      D='Fleißigen Bienen Bzz Bzz';    E=45479; F=3; G=4; H=0}
    @{A='Valdemar'; B='André';    C='';      D='Zuckerlager voll';            E=45479; F=1; G=3; H=0}
    @{A='Emilio';   B='Philipp';  C='Daisy'; D='Alles Andere als Arbeit';     E=45479; F=3; G=3; H=0}
    @{A='Chris';    B='Phia';     C='';      D='DD - Dirk & Dora';            E=45479; F=2; G=3; H=0}
    @{A='Phia';     B='André';    C='';      D='Team Heul doch';              E=45479; F=1; G=2; H=0}
    @{A='Valdemar'; B='Emilio';   C='';      D='Die romantischen Matrosen';   E=45479; F=2; G=2; H=0}
    @{A='Chris';    B='Phia';     C='';      D='ZaZa Grill';                  E=45500; F=2; G=3; H=0}
    @{A='Valdemar'; B='Emilio';   C='';      D='I got bit by a WIDDER';       E=45500; F=1; G=3; H=0}
    @{A='André';    B='Marian';   C='';      D='Bowle Batallion';             E=45500; F=3; G=3; H=0}
    @{A='Emilio';   B='Phia';     C='';      D='K-Hole';                      E=45500; F=1; G=4; H=0}
    @{A='André';    B='Leonie';   C='';      D='Schnelle Bälle';              E=45500; F=3; G=4; H=0}
    @{A='Chris';    B='Marian';   C='';      D='Best LoL-Players in the room';E=45500; F=3; G=4; H=0}
    @{A='Merlin';   B='Valdemar'; C='';      D='Two Bikey Boys Go Vroom';     E=45500; F=3; G=4; H=0}
)

$startRow = 160
$i = 0
foreach ($r in $rows) {
    $rowNum = $startRow + $i
    $ws.Cells.Item($rowNum, 1).Value() = $r.A
    $ws.Cells.Item($rowNum, 2).Value() = $r.B
    if ($r.C -ne '') {
        $ws.Cells.Item($rowNum, 3).Value() = $r.C
    }
    $ws.Cells.Item($rowNum, 4).Value() = $r.D
    $ws.Cells.Item($rowNum, 5).Value() = $r.E
    $ws.Cells.Item($rowNum, 6).Value() = $r.F
    $ws.Cells.Item($rowNum, 7).Value() = $r.G
    $ws.Cells.Item($rowNum, 8).Value() = $r.H
    $i = $i + 1
}

$lastRow = $startRow + $rows.Count - 1

# Match the date format (m/d/yyyy, same style as the rest of column E) by
# copying the format from an existing date cell onto the freshly written ones.
$ws.Cells.Item(2, 5).Copy()
$ws.Range("E160:E$lastRow").PasteSpecial(-4122)

# Move the selection to the first empty row below the new data, matching
# where the author's cursor ended up after entering the new rows.
$ws.Range("A176").Select()
